# Auto-generated edit script: updates profit-calculation columns (H-N)
# across multiple job sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, "H").Value = 445.2
$ws.Cells.Item(53, "I").Value = 389
$ws.Cells.Item(53, "J").Value = 529.5
$ws.Cells.Item(53, "K").Value = 389
$ws.Cells.Item(53, "L").Value = 529.5
$ws.Cells.Item(53, "M").Value = 248
$ws.Cells.Item(53, "N").Value = -1803.5
$ws.Cells.Item(106, "H").Value = 3101.5
$ws.Cells.Item(106, "I").Value = 3260
$ws.Cells.Item(106, "J").Value = 2837.3333
$ws.Cells.Item(106, "K").Value = 3260
$ws.Cells.Item(106, "L").Value = 2837.3333
$ws.Cells.Item(106, "M").Value = -2629
$ws.Cells.Item(106, "N").Value = -4099.3333
$ws.Cells.Item(111, "H").Value = 5414.84
$ws.Cells.Item(111, "I").Value = 6264.85
$ws.Cells.Item(111, "J").Value = 2014.8
$ws.Cells.Item(111, "K").Value = 18794.55
$ws.Cells.Item(111, "L").Value = 6044.4
$ws.Cells.Item(111, "M").Value = -15727.55
$ws.Cells.Item(111, "N").Value = -12178.4
$ws.Cells.Item(116, "H").Value = 2593.3333
$ws.Cells.Item(116, "I").Value = 1980
$ws.Cells.Item(116, "K").Value = 1980
$ws.Cells.Item(116, "M").Value = 1462
$ws.Cells.Item(127, "H").Value = 33334910
$ws.Cells.Item(127, "J").Value = 41668510
$ws.Cells.Item(127, "L").Value = 125005530
$ws.Cells.Item(127, "N").Value = -125015450
$ws.Cells.Item(129, "H").Value = 2397.9
$ws.Cells.Item(129, "I").Value = 5260.143
$ws.Cells.Item(129, "J").Value = 856.6923
$ws.Cells.Item(129, "K").Value = 15780.429
$ws.Cells.Item(129, "L").Value = 2570.0769
$ws.Cells.Item(129, "M").Value = -10780.429
$ws.Cells.Item(129, "N").Value = -12570.0769
$ws.Cells.Item(132, "H").Value = 4812252
$ws.Cells.Item(132, "I").Value = 5323872.5
$ws.Cells.Item(132, "J").Value = 3019.8
$ws.Cells.Item(132, "K").Value = 15971617.5
$ws.Cells.Item(132, "L").Value = 9059.400000000001
$ws.Cells.Item(132, "M").Value = -15969087.5
$ws.Cells.Item(132, "N").Value = -14119.4
$ws.Cells.Item(137, "H").Value = 1629.7826
$ws.Cells.Item(137, "I").Value = 1818.2667
$ws.Cells.Item(137, "J").Value = 1276.375
$ws.Cells.Item(137, "K").Value = 5454.800099999999
$ws.Cells.Item(137, "L").Value = 3829.125
$ws.Cells.Item(137, "M").Value = -2904.800099999999
$ws.Cells.Item(137, "N").Value = -8929.125
$ws.Cells.Item(138, "H").Value = 1775.909
$ws.Cells.Item(138, "I").Value = 1847.5333
$ws.Cells.Item(138, "J").Value = 1716.2222
$ws.Cells.Item(138, "K").Value = 5542.5999
$ws.Cells.Item(138, "L").Value = 5148.6666
$ws.Cells.Item(138, "M").Value = -402.5999000000002
$ws.Cells.Item(138, "N").Value = -15428.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, "H").Value = 7100.648
$ws.Cells.Item(32, "I").Value = 6103.6343
$ws.Cells.Item(32, "J").Value = 20726.5
$ws.Cells.Item(32, "K").Value = 6103.6343
$ws.Cells.Item(32, "L").Value = 20726.5
$ws.Cells.Item(32, "M").Value = -5816.6343
$ws.Cells.Item(32, "N").Value = -21300.5
$ws.Cells.Item(58, "H").Value = 13699.75
$ws.Cells.Item(58, "J").Value = 13699.75
$ws.Cells.Item(58, "L").Value = 13699.75
$ws.Cells.Item(58, "N").Value = -14559.75
$ws.Cells.Item(74, "H").Value = 1152.7273
$ws.Cells.Item(74, "I").Value = 1060
$ws.Cells.Item(74, "J").Value = 1400
$ws.Cells.Item(74, "K").Value = 1060
$ws.Cells.Item(74, "L").Value = 1400
$ws.Cells.Item(74, "M").Value = -186
$ws.Cells.Item(74, "N").Value = -3148
$ws.Cells.Item(77, "H").Value = 1152.7273
$ws.Cells.Item(77, "I").Value = 1060
$ws.Cells.Item(77, "J").Value = 1400
$ws.Cells.Item(77, "K").Value = 5300
$ws.Cells.Item(77, "L").Value = 7000
$ws.Cells.Item(77, "M").Value = -932
$ws.Cells.Item(77, "N").Value = -15736
$ws.Cells.Item(132, "H").Value = 2712.3386
$ws.Cells.Item(132, "I").Value = 2714.5625
$ws.Cells.Item(132, "J").Value = 2704.7144
$ws.Cells.Item(132, "K").Value = 8143.6875
$ws.Cells.Item(132, "L").Value = 8114.1432
$ws.Cells.Item(132, "M").Value = -5613.6875
$ws.Cells.Item(132, "N").Value = -13174.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, "H").Value = 70562
$ws.Cells.Item(86, "I").Value = 93342.664
$ws.Cells.Item(86, "K").Value = 93342.664
$ws.Cells.Item(86, "M").Value = -92219.664
$ws.Cells.Item(89, "H").Value = 70562
$ws.Cells.Item(89, "I").Value = 93342.664
$ws.Cells.Item(89, "K").Value = 466713.32
$ws.Cells.Item(89, "M").Value = -461097.32
$ws.Cells.Item(134, "H").Value = 1921.1459
$ws.Cells.Item(134, "I").Value = 1630.9318
$ws.Cells.Item(134, "J").Value = 5113.5
$ws.Cells.Item(134, "K").Value = 4892.7954
$ws.Cells.Item(134, "L").Value = 15340.5
$ws.Cells.Item(134, "M").Value = -2357.7954
$ws.Cells.Item(134, "N").Value = -20410.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, "H").Value = 133.61111
$ws.Cells.Item(7, "I").Value = 99.90909000000001
$ws.Cells.Item(7, "K").Value = 99.90909000000001
$ws.Cells.Item(7, "M").Value = 13.09090999999999
$ws.Cells.Item(134, "H").Value = 1393.8889
$ws.Cells.Item(134, "I").Value = 1401.4546
$ws.Cells.Item(134, "J").Value = 1382
$ws.Cells.Item(134, "K").Value = 4204.3638
$ws.Cells.Item(134, "L").Value = 4146
$ws.Cells.Item(134, "M").Value = -1669.3638
$ws.Cells.Item(134, "N").Value = -9216

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, "H").Value = 723.93335
$ws.Cells.Item(113, "J").Value = 513.6667
$ws.Cells.Item(113, "L").Value = 1541.0001
$ws.Cells.Item(113, "N").Value = -5881.0001
$ws.Cells.Item(131, "H").Value = 798.3711499999999
$ws.Cells.Item(131, "J").Value = 807.8936
$ws.Cells.Item(131, "L").Value = 2423.6808
$ws.Cells.Item(131, "N").Value = -12503.6808

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, "H").Value = 76925910
$ws.Cells.Item(97, "I").Value = 111113920
$ws.Cells.Item(97, "K").Value = 111113920
$ws.Cells.Item(97, "M").Value = -111113424
$ws.Cells.Item(126, "H").Value = 2246.6667
$ws.Cells.Item(126, "I").Value = 2353.923
$ws.Cells.Item(126, "J").Value = 2072.375
$ws.Cells.Item(126, "K").Value = 7061.768999999999
$ws.Cells.Item(126, "L").Value = 6217.125
$ws.Cells.Item(126, "M").Value = -4591.768999999999
$ws.Cells.Item(126, "N").Value = -11157.125
$ws.Cells.Item(132, "H").Value = 3137.4348
$ws.Cells.Item(132, "I").Value = 2565.5386
$ws.Cells.Item(132, "J").Value = 3880.9
$ws.Cells.Item(132, "K").Value = 7696.6158
$ws.Cells.Item(132, "L").Value = 11642.7
$ws.Cells.Item(132, "M").Value = -5166.6158
$ws.Cells.Item(132, "N").Value = -16702.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, "H").Value = 356140.4
$ws.Cells.Item(16, "I").Value = 67618.734
$ws.Cells.Item(16, "J").Value = 626629.5
$ws.Cells.Item(16, "K").Value = 67618.734
$ws.Cells.Item(16, "L").Value = 626629.5
$ws.Cells.Item(16, "M").Value = -67448.734
$ws.Cells.Item(16, "N").Value = -626969.5
$ws.Cells.Item(22, "H").Value = 2353.3333
$ws.Cells.Item(22, "J").Value = 5000
$ws.Cells.Item(22, "L").Value = 5000
$ws.Cells.Item(22, "N").Value = -5590
$ws.Cells.Item(27, "H").Value = 2353.3333
$ws.Cells.Item(27, "J").Value = 5000
$ws.Cells.Item(27, "L").Value = 5000
$ws.Cells.Item(27, "N").Value = -5214
$ws.Cells.Item(40, "H").Value = 62364.883
$ws.Cells.Item(40, "I").Value = 148003.28
$ws.Cells.Item(40, "J").Value = 2418
$ws.Cells.Item(40, "K").Value = 148003.28
$ws.Cells.Item(40, "L").Value = 2418
$ws.Cells.Item(40, "M").Value = -147867.28
$ws.Cells.Item(40, "N").Value = -2690
$ws.Cells.Item(46, "H").Value = 2530545
$ws.Cells.Item(46, "J").Value = 3373896.8
$ws.Cells.Item(46, "L").Value = 3373896.8
$ws.Cells.Item(46, "N").Value = -3374272.8
$ws.Cells.Item(132, "H").Value = 3898.3794
$ws.Cells.Item(132, "I").Value = 3749.25
$ws.Cells.Item(132, "K").Value = 11247.75
$ws.Cells.Item(132, "M").Value = -8717.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, "H").Value = 8769.299999999999
$ws.Cells.Item(15, "J").Value = 8921.444
$ws.Cells.Item(15, "L").Value = 8921.444
$ws.Cells.Item(15, "N").Value = -9497.444
$ws.Cells.Item(62, "H").Value = 6995933.5
$ws.Cells.Item(62, "I").Value = 25644992
$ws.Cells.Item(62, "J").Value = 2536.625
$ws.Cells.Item(62, "K").Value = 25644992
$ws.Cells.Item(62, "L").Value = 2536.625
$ws.Cells.Item(62, "M").Value = -25644368
$ws.Cells.Item(62, "N").Value = -3784.625
$ws.Cells.Item(65, "H").Value = 6995933.5
$ws.Cells.Item(65, "I").Value = 25644992
$ws.Cells.Item(65, "J").Value = 2536.625
$ws.Cells.Item(65, "K").Value = 128224960
$ws.Cells.Item(65, "L").Value = 12683.125
$ws.Cells.Item(65, "M").Value = -128221840
$ws.Cells.Item(65, "N").Value = -18923.125
$ws.Cells.Item(122, "H").Value = 2394.8
$ws.Cells.Item(122, "I").Value = 1401
$ws.Cells.Item(122, "J").Value = 3057.3333
$ws.Cells.Item(122, "K").Value = 4203
$ws.Cells.Item(122, "L").Value = 9171.999899999999
$ws.Cells.Item(122, "M").Value = -1753
$ws.Cells.Item(122, "N").Value = -14071.9999
$ws.Cells.Item(123, "H").Value = 30000.5
$ws.Cells.Item(123, "J").Value = 30000.5
$ws.Cells.Item(123, "L").Value = 30000.5
$ws.Cells.Item(123, "N").Value = -39800.5
$ws.Cells.Item(126, "H").Value = 1468.7142
$ws.Cells.Item(126, "I").Value = 1235.375
$ws.Cells.Item(126, "J").Value = 1779.8334
$ws.Cells.Item(126, "K").Value = 3706.125
$ws.Cells.Item(126, "L").Value = 5339.5002
$ws.Cells.Item(126, "M").Value = -1236.125
$ws.Cells.Item(126, "N").Value = -10279.5002

